# Refresh the cryptos price/volume snapshot (columns D and E) to match
# the latest scrape. Only the cells that actually changed are touched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.994.26'
$ws.Range('E2').Value = '  -0.77%  '
# Row 3
$ws.Range('D3').Value = '2.551.67'
# Row 4
$ws.Range('E4').Value = '  -0.01%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.92'
$ws.Range('E5').Value = '  +1.38%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.13'
$ws.Range('E6').Value = '  -1.81%  '
# Row 7
$ws.Range('E7').Value = '  -0.02%  '
# Row 8
$ws.Range('E8').Value = '  -0.79%  '
# Row 9
$ws.Range('E9').Value = '  -0.98%  '
# Row 10
$ws.Range('E10').Value = '  -4.29%  '
# Row 11
$ws.Range('E11').Value = '  -0.30%  '
# Row 12
$ws.Range('E12').Value = '  -1.38%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.20'
$ws.Range('E13').Value = '  -3.10%  '
# Row 14
$ws.Range('D14').Value = '3.008.97'
$ws.Range('E14').Value = '  -0.09%  '
# Row 15
$ws.Range('D15').Value = '62.912.12'
$ws.Range('E15').Value = '  -0.85%  '
# Row 16
$ws.Range('E16').Value = '  -0.70%  '
# Row 17
$ws.Range('D17').Value = '2.561.90'
$ws.Range('E17').Value = '  -0.29%  '
# Row 18
$ws.Range('E18').Value = '  -2.44%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '335.68'
$ws.Range('E19').Value = '  -1.73%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.33'
$ws.Range('E20').Value = '  -0.82%  '
# Row 21
$ws.Range('E21').Value = '  -2.24%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.15%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.53'
$ws.Range('E23').Value = '  -1.04%  '
# Row 24
$ws.Range('E24').Value = '  -0.43%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.61'
$ws.Range('E25').Value = '  +2.43%  '
# Row 26
$ws.Range('E26').Value = '  +0.95%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.997'
$ws.Range('E27').Value = '  -0.43%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.39'
$ws.Range('E28').Value = '  -0.55%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.32'
$ws.Range('E29').Value = '  +1.29%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.89'
$ws.Range('E30').Value = '  +1.17%  '
# Row 31
$ws.Range('E31').Value = '  -3.13%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '177.18'
$ws.Range('E32').Value = '  -0.24%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.55'
$ws.Range('E33').Value = '  -1.84%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '410.60'
$ws.Range('E34').Value = '  -0.90%  '
# Row 35
$ws.Range('E35').Value = '  +0.24%  '
# Row 36
$ws.Range('E36').Value = '  -1.02%  '
# Row 38
$ws.Range('E38').Value = '  -2.66%  '
# Row 39
$ws.Range('E39').Value = '  -0.45%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.11%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '39.68'
$ws.Range('E41').Value = '  -0.78%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '151.18'
$ws.Range('E42').Value = '  -2.60%  '
# Row 43
$ws.Range('E43').Value = '  -1.34%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.80'
$ws.Range('E44').Value = '  -1.35%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0538'
$ws.Range('E45').Value = '  +0.81%  '
# Row 47
$ws.Range('E47').Value = '  -0.03%  '
# Row 48
$ws.Range('E48').Value = '  +1.93%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.25'
$ws.Range('E49').Value = '  -2.69%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.72'
$ws.Range('E50').Value = '  -7.80%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.30'
$ws.Range('E51').Value = '  -0.13%  '
